# Add two new columns to the right of the existing table:
#   I: header "I0", every data row = 1
#   J: header "IF", every data row = copy of column H's value
#
# Column H's header cell (H1) is used as the formatting template for the
# new header cells so they pick up the same style (bold, bordered,
# centered) as the rest of row 1, instead of the default style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -------------------------------------------------------------
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# --- Data rows 2-20 --------------------------------------------------------
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 9).Value2 = 1
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
